# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, and fixes the Mantle / VeChain row ordering (rows 48-49
# had swapped places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    # Assign as literal text. A leading apostrophe forces Excel to keep
    # numeric-looking strings (e.g. "1.00", "0.0000197") as text instead
    # of coercing them to a Double and dropping significant digits.
    $ws.Range($cellRef).Value = "'" + $value
}

# Row 2
$ws.Range("D2").Value = '65.674.42'
$ws.Range("E2").Value = '  -0.38%  '
# Row 3
$ws.Range("D3").Value = '2.671.39'
$ws.Range("E3").Value = '  +0.23%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
Set-Text "D5" '601.93'
$ws.Range("E5").Value = '  -1.18%  '
# Row 6
Set-Text "D6" '157.44'
$ws.Range("E6").Value = '  -0.35%  '
# Row 7
Set-Text "D7" '0.999'
$ws.Range("E7").Value = '  +0.00%  '
# Row 8
Set-Text "D8" '0.623'
$ws.Range("E8").Value = '  +5.58%  '
# Row 9
Set-Text "D9" '0.126'
$ws.Range("E9").Value = '  +0.42%  '
# Row 10
Set-Text "D10" '0.402'
$ws.Range("E10").Value = '  -0.56%  '
# Row 11
$ws.Range("E11").Value = '  -3.32%  '
# Row 12
Set-Text "D12" '0.155'
$ws.Range("E12").Value = '  -0.22%  '
# Row 13
Set-Text "D13" '29.45'
$ws.Range("E13").Value = '  -2.24%  '
# Row 14
Set-Text "D14" '0.0000197'
$ws.Range("E14").Value = '  -5.96%  '
# Row 15
$ws.Range("D15").Value = '3.150.02'
$ws.Range("E15").Value = '  +0.15%  '
# Row 16
$ws.Range("D16").Value = '65.496.35'
$ws.Range("E16").Value = '  -0.38%  '
# Row 17
$ws.Range("D17").Value = '2.674.38'
$ws.Range("E17").Value = '  +0.40%  '
# Row 18
Set-Text "D18" '12.71'
$ws.Range("E18").Value = '  +0.21%  '
# Row 19
Set-Text "D19" '4.82'
$ws.Range("E19").Value = '  -1.55%  '
# Row 20
Set-Text "D20" '7.69'
$ws.Range("E20").Value = '  +3.23%  '
# Row 21
Set-Text "D21" '351.20'
$ws.Range("E21").Value = '  -2.34%  '
# Row 22
Set-Text "D22" '1.00'
$ws.Range("E22").Value = '  -0.06%  '
# Row 23
Set-Text "D23" '69.47'
$ws.Range("E23").Value = '  -1.10%  '
# Row 24
Set-Text "D24" '0.0000110'
$ws.Range("E24").Value = '  +2.62%  '
# Row 25
Set-Text "D25" '9.77'
$ws.Range("E25").Value = '  +2.53%  '
# Row 26
Set-Text "D26" '1.63'
$ws.Range("E26").Value = '  -3.85%  '
# Row 27
$ws.Range("E27").Value = '  -2.84%  '
# Row 28
Set-Text "D28" '1.59'
$ws.Range("E28").Value = '  -3.33%  '
# Row 29
Set-Text "D29" '8.07'
$ws.Range("E29").Value = '  -0.61%  '
# Row 30
$ws.Range("E30").Value = '  +0.33%  '
# Row 31
Set-Text "D31" '531.09'
$ws.Range("E31").Value = '  -1.18%  '
# Row 32
Set-Text "D32" '2.14'
$ws.Range("E32").Value = '  -2.57%  '
# Row 33
Set-Text "D33" '1.76'
$ws.Range("E33").Value = '  -2.01%  '
# Row 34
Set-Text "D34" '6.55'
$ws.Range("E34").Value = '  +1.20%  '
# Row 35
Set-Text "D35" '5.49'
$ws.Range("E35").Value = '  -0.81%  '
# Row 36
Set-Text "D36" '0.423'
$ws.Range("E36").Value = '  -2.46%  '
# Row 37
Set-Text "D37" '20.42'
$ws.Range("E37").Value = '  -1.65%  '
# Row 38
$ws.Range("E38").Value = '  +0.07%  '
# Row 39
Set-Text "D39" '159.27'
$ws.Range("E39").Value = '  -1.98%  '
# Row 40
Set-Text "D40" '1.94'
$ws.Range("E40").Value = '  -3.99%  '
# Row 41
Set-Text "D41" '1.00'
$ws.Range("E41").Value = '  +0.05%  '
# Row 42
Set-Text "D42" '42.79'
$ws.Range("E42").Value = '  +1.39%  '
# Row 43
Set-Text "D43" '165.27'
$ws.Range("E43").Value = '  -2.54%  '
# Row 44
Set-Text "D44" '4.10'
$ws.Range("E44").Value = '  -1.79%  '
# Row 45
Set-Text "D45" '0.0611'
$ws.Range("E45").Value = '  -0.21%  '
# Row 46
Set-Text "D46" '2.30'
$ws.Range("E46").Value = '  -2.69%  '
# Row 47
Set-Text "D47" '23.04'
$ws.Range("E47").Value = '  -0.30%  '
# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-Text "D48" '0.0260'
$ws.Range("E48").Value = '  -1.80%  '
# Row 49
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-Text "D49" '0.644'
$ws.Range("E49").Value = '  -2.57%  '
# Row 50
$ws.Range("E50").Value = '  +3.47%  '
# Row 51
Set-Text "D51" '20.26'
$ws.Range("E51").Value = '  +1.79%  '
